# Entrega 2.docx — "4.5 ESTIMATIVAS DE TAMANHO E ESFORÇO" table:
# fill in the previously-empty "Responsavel" / "Status" cells for the
# rows that belong to that section.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: "4.5 ESTIMATIVAS DE TAMANHO E ESFORÇO"
$t.Cell(1, 2).Range.Text = "Vinicius"
$t.Cell(1, 3).Range.Text = "Fazendo"

# Row 3: "6.1.1 Diagrama de Pacotes"
$t.Cell(3, 2).Range.Text = "Lucas"
$t.Cell(3, 3).Range.Text = "Fazendo"

# Row 4: "6.1.2 Diagramas de Classes"
$t.Cell(4, 2).Range.Text = "Lucas"

# Row 5: "6.1.3 Diagramas de Objetos"
$t.Cell(5, 2).Range.Text = "Lucas"
